$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "联赛名称"
$ws.Range("B1").Value = "联赛Id"
$ws.Range("C1").Value = "赛季"
$ws.Range("D1").Value = "比赛日"
$ws.Range("E1").Value = "比赛时间"
